# Convert the "m:comment" field (fldChar begin/instrText.../fldChar end)
# in the second paragraph into plain literal text runs, e.g.
#   { m : <bookmark> comment some important comment }
# while keeping the _GoBack bookmark that sits between the ":" and the
# leading space of the comment text.

$d = $word.ActiveDocument

# Locate the paragraph that holds the field (the one whose Fields.Count
# for its range is 1) instead of hard-coding an index, so the script is
# resilient to the exact paragraph numbering.
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Fields.Count -gt 0) {
        $targetPara = $para
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not find the paragraph containing the comment field."
}

$start = $targetPara.Range.Start
$end = $targetPara.Range.End

# Exclude the trailing paragraph mark so the <w:p> element (and its
# attributes) are preserved and only its children are replaced.
$rng = $d.Range($start, $end - 1)

$inner = '<w:r><w:t>{</w:t></w:r>' + `
         '<w:r><w:t>m</w:t></w:r>' + `
         '<w:r><w:t>:</w:t></w:r>' + `
         '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
         '<w:bookmarkEnd w:id="0"/>' + `
         '<w:r><w:t xml:space="preserve"> </w:t></w:r>' + `
         '<w:r><w:t>comment some important comment</w:t></w:r>' + `
         '<w:r><w:t xml:space="preserve">}</w:t></w:r>'

$snippet = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
           '<pkg:xmlData>' + `
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:body><w:p w:rsidP="00F5495F" w:rsidR="00C52979" w:rsidRDefault="00C52979">' + $inner + '</w:p></w:body>' + `
           '</w:document></pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($snippet) | Out-Null
